# Leave Card update: a new "SL(1-0-0)" leave entry (1 day, dated 3/8/2023)
# is inserted into Table1 immediately above the former row 79
# (PERIOD 4/1/2023), pushing all subsequent table rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# Grow the table boundary by one row *before* shifting cells so the
# structured-reference formulas already inside the table (e.g. the
# last row's EARNED-column IF/ISBLANK formula) are left untouched.
$tbl.Resize($ws.Range("A8:K133"))

# Shift row 79 (and everything below it) down by one row.
$ws.Rows.Item(79).Insert()

# The freshly inserted row 79 has no formatting yet - clone it from the
# row directly below (the row that used to be row 79, now row 80), which
# carries the standard "blank data row" styling used throughout the table.
$ws.Range("A80:K80").Copy()
$ws.Range("A79:K79").PasteSpecial(-4122)

# Column K needs the date-formatted style instead - copy that from K78.
$ws.Range("K78").Copy()
$ws.Range("K79").PasteSpecial(-4122)

# Fill in the new leave-record row's data.
$ws.Range("B79").Value = "SL(1-0-0)"
$ws.Range("H79").Value = 1
$ws.Range("K79").Value = 44993
$ws.Range("G79").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Restore the active-cell selection to where it now sits one row lower.
$null = $ws.Range("B80").Select()
